$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old rows (2-8) that held the "Location"/"Year" themed Q&A data
$ws.Range("A2:C8").ClearContents()

# Write the new "Person" themed Q&A rows
$data = @(
    @("Who was the F1 World Champion in 2022?", "Max Verstappen", "Person"),
    @("Who is the artist behind the song Thriller?", "Michael Jackson", "Person"),
    @("Who is the current president of the US?", "Joe Biden", "Person")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("L9").Select() | Out-Null
